$p = $ppt.ActivePresentation

# Add a new slide at the end (position 5) using the "Title and Content"
# layout (ppLayoutText = 2) -- the same layout used by slides 2-4.
$s = $p.Slides.Add(5, 2)

# --- Title placeholder: "Task Creation", centered ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Task Creation"
$title.ParagraphFormat.Alignment = 2

# --- Body / content placeholder ---
$body = $s.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: "Task Scaffolding – Pre-" | "reqs" | " and setup for task creation"
$body.Text = "Task Scaffolding – Pre-"
$body.InsertAfter("reqs") | Out-Null
$body.InsertAfter(" and setup for task creation") | Out-Null

# Paragraph 2
$body.InsertAfter("`rCreate Task") | Out-Null

# Paragraph 3
$body.InsertAfter("`rPackage Task") | Out-Null

# Paragraph 4
$body.InsertAfter("`rPublish Task") | Out-Null

# Paragraph 5: "Sharing Task " | "with Organizations"
$body.InsertAfter("`rSharing Task ") | Out-Null
$body.InsertAfter("with Organizations") | Out-Null

Write-Output "Added slide 5 (Task Creation)"
